# loginTestDataPW.xlsx - update to docker/yml files
# - validUsers!A2/B2 get a new valid login (email + password)
# - validUsers becomes the active sheet/tab, with B3 selected

$wb = $excel.ActiveWorkbook

$validUsers = $wb.Worksheets.Item("validUsers")
$invalidUsers = $wb.Worksheets.Item("invalidUsers")

# Update the valid-user credentials used by the tests
$validUsers.Range("A2").Value = "natashatestpw@gmail.com"
$validUsers.Range("B2").Value = "3230474N5a5t5e5!"

# Make sure invalidUsers keeps its previous selection before we switch away
$invalidUsers.Range("F11").Select()

# validUsers is now the active/selected sheet, with B3 selected
$validUsers.Activate()
$validUsers.Range("B3").Select()
